# Add a new "2022-Q1" sheet (fund holdings) positioned between the existing
# "2021-Q2" sheet and the "总计" (totals) sheet, and add a corresponding
# summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" worksheet right before "总计" -------------
$totalSheetBefore = $wb.Worksheets("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: sheet handles here resolve by position, so after the insert above
# changed sheet ordering/count we must re-fetch "总计" by name rather than
# reuse $totalSheetBefore (which now resolves to the new sheet's slot).
$totalSheet = $wb.Worksheets("总计")

# --- 2. Header row (copy header style from the "总计" sheet's header) ------
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- 3. Data rows ------------------------------------------------------------
# Style A2:A3 like the "总计" sheet's A2 (row-index marker column)
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# Text-valued columns must stay text (preserve leading zeros / exact digits)
$textRange = $newSheet.Range("B2:G3")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "014307"
$newSheet.Range("C2").Value = "嘉实多元动力混合A"
$newSheet.Range("D2").Value = "1.83"
$newSheet.Range("E2").Value = "91.81"
$newSheet.Range("F2").Value = "4.20"
$newSheet.Range("G2").Value = "0.0769"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "014308"
$newSheet.Range("C3").Value = "嘉实多元动力混合C"
$newSheet.Range("D3").Value = "0.20"
$newSheet.Range("E3").Value = "91.81"
$newSheet.Range("F3").Value = "4.20"
$newSheet.Range("G3").Value = "0.0084"
$newSheet.Range("H3").Value = 10

# --- 4. Update "总计": insert a new row above the existing "2021-Q2" row ----
$totalSheet.Rows.Item(2).Insert(-4121)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.09

# the row-index marker in column A is a 0-based row counter; the existing
# "2021-Q2" row (now shifted down to row 3) moves from index 0 to index 1
$totalSheet.Range("A3").Value = 1
